$d = $word.ActiveDocument

$d.Content.Find.Execute("44×17=", $true, $false, $false, $false, $false, $true, 1, $false, "37×73=", 2) | Out-Null
$d.Content.Find.Execute("11×94=", $true, $false, $false, $false, $false, $true, 1, $false, "74×51=", 2) | Out-Null
$d.Content.Find.Execute("72×48=", $true, $false, $false, $false, $false, $true, 1, $false, "30×32=", 2) | Out-Null
$d.Content.Find.Execute("12×36=", $true, $false, $false, $false, $false, $true, 1, $false, "78×33=", 2) | Out-Null
$d.Content.Find.Execute("86×34=", $true, $false, $false, $false, $false, $true, 1, $false, "26×44=", 2) | Out-Null
$d.Content.Find.Execute("52×68=", $true, $false, $false, $false, $false, $true, 1, $false, "83×24=", 2) | Out-Null
$d.Content.Find.Execute("12×66=", $true, $false, $false, $false, $false, $true, 1, $false, "98×85=", 2) | Out-Null
$d.Content.Find.Execute("12×48=", $true, $false, $false, $false, $false, $true, 1, $false, "89×50=", 2) | Out-Null
$d.Content.Find.Execute("82×23=", $true, $false, $false, $false, $false, $true, 1, $false, "11×20=", 2) | Out-Null
$d.Content.Find.Execute("31×96=", $true, $false, $false, $false, $false, $true, 1, $false, "80×96=", 2) | Out-Null
$d.Content.Find.Execute("56×29=", $true, $false, $false, $false, $false, $true, 1, $false, "57×56=", 2) | Out-Null
$d.Content.Find.Execute("72×64=", $true, $false, $false, $false, $false, $true, 1, $false, "34×30=", 2) | Out-Null
$d.Content.Find.Execute("46×59=", $true, $false, $false, $false, $false, $true, 1, $false, "73×58=", 2) | Out-Null
$d.Content.Find.Execute("80×60=", $true, $false, $false, $false, $false, $true, 1, $false, "34×36=", 2) | Out-Null
$d.Content.Find.Execute("95×89=", $true, $false, $false, $false, $false, $true, 1, $false, "84×28=", 2) | Out-Null
$d.Content.Find.Execute("74×80=", $true, $false, $false, $false, $false, $true, 1, $false, "12×86=", 2) | Out-Null
$d.Content.Find.Execute("65×20=", $true, $false, $false, $false, $false, $true, 1, $false, "23×30=", 2) | Out-Null
$d.Content.Find.Execute("74×65=", $true, $false, $false, $false, $false, $true, 1, $false, "68×39=", 2) | Out-Null
$d.Content.Find.Execute("24×49=", $true, $false, $false, $false, $false, $true, 1, $false, "65×25=", 2) | Out-Null
$d.Content.Find.Execute("27×21=", $true, $false, $false, $false, $false, $true, 1, $false, "99×70=", 2) | Out-Null
$d.Content.Find.Execute("78×45=", $true, $false, $false, $false, $false, $true, 1, $false, "18×93=", 2) | Out-Null
$d.Content.Find.Execute("17×36=", $true, $false, $false, $false, $false, $true, 1, $false, "55×86=", 2) | Out-Null
$d.Content.Find.Execute("81×47=", $true, $false, $false, $false, $false, $true, 1, $false, "63×30=", 2) | Out-Null
$d.Content.Find.Execute("86×17=", $true, $false, $false, $false, $false, $true, 1, $false, "40×87=", 2) | Out-Null
$d.Content.Find.Execute("23×51=", $true, $false, $false, $false, $false, $true, 1, $false, "12×91=", 2) | Out-Null
